# Insert a new data row into the "Zapallo italiano" price series.
# A new record is inserted before the existing row 243, shifting the
# existing rows 243:284 down to 244:285 (dimension grows from R284 to R285).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 243:284 down by one row, creating a blank row 243.
$ws.Rows.Item(243).EntireRow.Insert()

# Populate the newly inserted row 243 with the new record. Columns not
# listed below (A,B,C,E,F,G,H,I,J,N,O,Q,R) carry the same values the old
# row 243 had, since Insert() already copied that row's formatting/values
# down - so only the cells that actually differ need to be (re)written here,
# but to be explicit and safe we set every cell in the row.
$ws.Range("A243").Value = 11
$ws.Range("B243").Value = "Vega Monumental Concepción"
$ws.Range("C243").Value = "Bíobío"
$ws.Range("D243").Value = 45218
$ws.Range("D243").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E243").Value = 8
$ws.Range("F243").Value = 100112032
$ws.Range("G243").Value = "Zapallo italiano"
$ws.Range("H243").Value = "Sin especificar"
$ws.Range("I243").Value = "Primera"
$ws.Range("J243").Value = 100
$ws.Range("K243").Value = 17000
$ws.Range("L243").Value = 18000
$ws.Range("M243").Value = 17500
$ws.Range("N243").Value = "`$/caja 50 unidades"
$ws.Range("O243").Value = "Región de Arica y Parinacota"
$ws.Range("P243").Value = 350
$ws.Range("Q243").Value = 50
$ws.Range("R243").Value = "Hortaliza"
